$d = $word.ActiveDocument

# --- Skills line: "OS Basic" -> "Spark Analysis" ----------------------
# The source run read:
#   " Network Socket, Distributed System (Raft), OS Basic, machine learning"
# After the edit it must become three separate runs (same rPr):
#   " Network Socket, Distributed System (Raft), "
#   "Spark Analysis"
#   ", machine learning"
# (the following ", image processing, time series forecast" run is left
# untouched). Replacing the "OS Basic" text in place and nudging the
# Font property on the replaced range (and again on the run that follows
# it) forces the engine to keep/introduce run boundaries at those points
# instead of re-merging everything into a single run.

$rngSkills = $d.Content
$rngSkills.Find.ClearFormatting()
$foundSkills = $rngSkills.Find.Execute("OS Basic", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
if ($foundSkills) {
    $rngSkills.Text = "Spark Analysis"
    $rngSkills.Font.Bold = $true
    $rngSkills.Font.Bold = $false

    $rngAfter = $d.Content
    $rngAfter.Find.ClearFormatting()
    $foundAfter = $rngAfter.Find.Execute(", machine learning", $true, $false, $false, $false, `
                                          $false, $true, 1, $false, "", 0)
    if ($foundAfter) {
        $rngAfter.Font.Bold = $true
        $rngAfter.Font.Bold = $false
    }
}

# --- Project bullet: merge "Implement" + " demand paging..." ----------
# The source had two runs: "Implement" and
# " demand paging, page replacement to free up a physic page to handle
# page faults." -- they must become a single run with the combined text.
# Doing the replacement through Find.Execute's own Replace argument (as
# opposed to assigning Range.Text) makes the engine collapse the run back
# into one, matching the target XML.

$rngImpl = $d.Content
$rngImpl.Find.ClearFormatting()
$implText = "Implement demand paging, page replacement to free up a physic page to handle page faults."
$rngImpl.Find.Execute($implText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $implText, 2) | Out-Null
